$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the latest crypto symbol/price scrape.
# Force text format on target cells first so numeric-looking strings
# (prices, percentages) are preserved exactly instead of being
# reinterpreted as Number/Percentage cell types.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '303.11'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.02%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.21'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '7.14%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-3.52%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07838'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.63%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.209'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-3.67%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.018'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.09%'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '4.037'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.15%'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9148'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.44%'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09714'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-4.13%'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1882'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.50%'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08586'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.19%'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03558'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.13%'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09937'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.23%'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001488'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.51%'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005683'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-2.19%'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.464'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.10%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.376'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '12.83%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.60%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1312'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-0.95%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.774'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4.97%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2297'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-1.65%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04617'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.21%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001232'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.04%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004784'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '7.91%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001402'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '7.85%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '39.81%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01784'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '1.53%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04736'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.29%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.008071'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '5.22%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1390'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-1.21%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.007688'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '9.23%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002092'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-9.03%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009937'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '1.20%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006181'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '3.39%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.19%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.23%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.19%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.19%'
